# regen sval data to filter save games
# Update the numeric stat columns (B:G) for rows 2-5 with the regenerated
# values that result from filtering out save games from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 9.983522426115931, 0.1496068669990043, 13.86384647080068, 0, 24.65154103579788)
    3 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    4 = @(0.2881169905109251, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1, 1.012145535086602)
    5 = @(0.1169995834814548, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 16.32892827181126)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
